$p = $ppt.ActivePresentation

function Get-ShapeByName($shapes, $name) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.Name -eq $name) {
            return $candidate
        }
    }
    return $shapes.Item(1)
}

# --- 1) Slide master footer date field: "19.06.2024" -> "06.07.2024" ---
$m = $p.SlideMaster
$dateShape = Get-ShapeByName $m.Shapes "Rectangle 6"
$dateShape.TextFrame.TextRange.Text = "06.07.2024"

# --- 2) Slide 1 subtitle: "19.06.2024, Daniel Krämer"
#        -> split into "10.07.2024" + ", Daniel Krämer" ---
$s = $p.Slides.Item(1)
$subtitle = Get-ShapeByName $s.Shapes "Rectangle 3"
$tr = $subtitle.TextFrame.TextRange
$datePart = $tr.Characters(1, 10)
$datePart.Text = "10.07.2024"
